# Apply the refreshed crypto price/volume snapshot (Sat Jan 14 09:30:02 UTC 2023 run).
# Only the "Price" (column D) and "Volume(1h)" (column E) cells for the affected rows
# change; everything else (Coin, Link, Data, Hora columns, row count, etc.) stays put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parallel arrays: cell reference -> new text value.
$cellRefs = @(
    'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47'
)
$newValues = @(
    '307.46', '7.08%', '32.10', '8.43%', '5.278', '2.94%', '0.07398', '10.38%', '7.907', '7.89%', '3.796', '11.56%', '1.497', '9.77%', '0.9122', '-0.74%', '0.01659', '2,470.75%', '0.1688', '6.06%', '0.07559', '11.77%', '0.07964', '3.14%', '0.03079', '5.24%', '0.09848', '9.61%', '0.001514', '-4.64%', '0.04559', '1.71%', '0.006268', '-0.35%', '3.467', '0.38%', '2.235', '0.40%', '0.3288', '2.27%', '0.1352', '3.79%', '4.303', '5.96%', '0.1648', '4.27%', '0.001229', '3.12%', '0.004462', '8.29%', '0.0001308', '9.22%', '0.0001753', '8.51%', '0.04519', '5.78%', '0.007074', '5.29%', '0.1351', '8.99%', '0.002275', '2.17%', '0.01405', '17.12%', '0.00006090', '7.43%', '-3.82%', '0.01309', '0.31%'
)

for ($i = 0; $i -lt $cellRefs.Length; $i++) {
    $cell = $ws.Range($cellRefs[$i])
    # Force text format before assigning so numeric-looking strings (e.g. "307.46",
    # "7.08%", "2,470.75%") are stored as literal text, not auto-converted to a
    # Number/Percentage value by Excel's input parser.
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$i]
    # Reset the style back to the workbook default so we do not leave a stray
    # "Text"/quote-prefix style on the cell (the source cells carry no explicit style).
    $cell.Style = "Normal"
}

